$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restructure columns -------------------------------------------------
# Before: A..F unchanged, G="2012-2021" (dup range), H="2021", I="2022 (projection)"
# After : A..F unchanged, G="2022", H="2023 (projection)"
# Drop the old "2022 (projection)" column (I) and the duplicate "2012-2021"
# column (G); this shifts the old "2021" column (H) left into G.
$ws.Columns("I").Delete()
$ws.Columns("G").Delete()

# --- Header row -----------------------------------------------------------
# G1 currently holds the old "2021" label; replace with "2022". Because
# "2022" looks like a pure number, a plain .Value assignment would store it
# as a numeric cell instead of text, so build it as text via a formula and
# paste the computed value back in (keeps it a real text cell, no leftover
# number formatting on the cell itself).
$ws.Range("Z1").Formula = "=TEXT(2022,""0"")"
$ws.Range("Z1").Copy()
$ws.Range("G1").PasteSpecial(-4163)
$ws.Range("H1").Value = "2023 (projection)"

# --- Row 2: GHG -------------------------------------------------------
$ws.Range("E2").Value = "46±5.2"
$ws.Range("F2").Value = "53±5.5"
$ws.Range("G2").Value = "55±5.4"

# --- Row 3: CO2-FFI -----------------------------------------------------
$ws.Range("G3").Value = "37±3"
$ws.Range("H3").Value = "X"

# --- Row 4: CO2-LULUCF ---------------------------------------------------
$ws.Range("G4").Value = "4.3±3"
$ws.Range("H4").Value = "X"

# --- Row 5: CH4 -----------------------------------------------------------
$ws.Range("D5").Value = "7.2±2.2"
$ws.Range("E5").Value = "8.1±2.4"
$ws.Range("F5").Value = "8.8±2.6"
$ws.Range("G5").Value = "9±2.7"

# --- Row 6: N2O -------------------------------------------------------
$ws.Range("D6").Value = "2.3±1.4"
$ws.Range("E6").Value = "2.7±1.6"
$ws.Range("F6").Value = "2.9±1.8"
$ws.Range("G6").Value = "3.1±1.9"

# --- Row 7: F-gases ---------------------------------------------------
$ws.Range("B7").Value = "0.53±0.16"
$ws.Range("C7").Value = "0.64±0.19"
$ws.Range("D7").Value = "0.55±0.16"
$ws.Range("E7").Value = "0.74±0.22"
$ws.Range("F7").Value = "1.1±0.34"
$ws.Range("G7").Value = "1.4±0.43"

# --- cleanup scratch cell --------------------------------------------------
$ws.Columns("Z").Delete()
